$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Global_variable": replace the multi-borehole (Hsinchu) rows with a
# single row describing the new "A2 (2018)" / JDN project, and clear out the
# other rows that used to hold the extra boreholes (BH01..BH08).
# ---------------------------------------------------------------------------
$gv = $wb.Worksheets.Item("Global_variable")

# Clear rows 3:9 (contents only, keep formatting) - these used to hold the
# extra boreholes BH03, BH04, BH02, BH05, BH06, BH07, BH08.
$gv.Range("A3:V9").ClearContents()

# Update row 2 with the new project values.
$gv.Range("B2").Value = "C:\Users\sc\PycharmProjects\SI_processing_automation_SC"
$gv.Range("C2").Value = "C:\Users\sc\PycharmProjects\SI_processing_automation_SC\Input_files\cpt_data_files\JDN"
$gv.Range("D2").Value = "C:\Users\sc\PycharmProjects\SI_processing_automation_SC\Output_files\CPT-fig\JDN"
$gv.Range("E2").Value = "C:\Users\sc\PycharmProjects\SI_processing_automation_SC\Output_files\Processed-CPT\JDN"
$gv.Range("G2").Value = $true
$gv.Range("H2").Value = $true
$gv.Range("I2").Value = $true
$gv.Range("L2").Value = "SCPG-combined.xlsx"
$gv.Range("N2").Value = "SOIL_PROPERTY-combined.xlsx"
$gv.Range("O2").Value = "A2 (2018)_CPT_processed_data.csv"
$gv.Range("R2").Value = "A2 (2018)"

# S2/T2 no longer hold extra SCPT locations.
$gv.Range("S2").ClearContents()
$gv.Range("T2").ClearContents()

$gv.Range("D1").Select()
$gv.Range("C8").Select()

# ---------------------------------------------------------------------------
# Sheet "Stratigraphy_color_dict": append a new "Dummy" unit reusing the
# "Silt" colour.
# ---------------------------------------------------------------------------
$sc = $wb.Worksheets.Item("Stratigraphy_color_dict")
$sc.Range("A10").Value = "Dummy"
$sc.Range("B10").Value = "#b0ebff"
$sc.Range("A9").Copy() | Out-Null
$sc.Range("A10").PasteSpecial(-4122) | Out-Null
$sc.Range("B10").Select()
